# Adding Excel file to Read data
# Mirrors the authors edit: company name + phone text fixed up, phone
# stored as a left/top-aligned text value, a few column widths tweaked,
# the active selection moved to D1 and the sheet's page orientation set
# to portrait.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell content -----------------------------------------------------
# B1: "prosthatic ltd" -> "Prosthatic pvt ltd"
$ws.Range("B1").Value = "Prosthatic pvt ltd"

# G1: "98765+21" -> "987654+908"
$ws.Range("G1").Value = "987654+908"

# D1: numeric phone number -> text value (quoted string), left/top
# aligned, formatted with the integer "0" number format.
$ws.Range("D1").NumberFormat = "0"
$ws.Range("D1").HorizontalAlignment = -4131
$ws.Range("D1").VerticalAlignment = -4160
$ws.Range("D1").Value = """9087654321"""

# --- Column widths ------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 20.166666666666668
$ws.Columns.Item(3).ColumnWidth = 17.666666666666668
$ws.Columns.Item(4).ColumnWidth = 16
$ws.Columns.Item(5).ColumnWidth = 11.166666666666666
$ws.Columns.Item(6).ColumnWidth = 20.833333333333332
$ws.Columns.Item(7).ColumnWidth = 11.833333333333334

# --- Selection ----------------------------------------------------------
$null = $ws.Range("D1").Select()

# --- Page setup -----------------------------------------------------------
$ws.PageSetup.Orientation = 1
